$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 140 (pushes existing rows 140.. down to 142..)
$ws.Rows.Item(140).Insert()
$ws.Rows.Item(140).Insert()

# New row 140 data
$ws.Range("A140").Value = 6
$ws.Range("B140").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C140").Value = "Metropolitana"
$ws.Range("D140").Value = 44505
$ws.Range("E140").Value = 13
$ws.Range("F140").Value = 100112022
$ws.Range("G140").Value = "Arveja Verde"
$ws.Range("H140").Value = "Perfection"
$ws.Range("I140").Value = "Primera"
$ws.Range("J140").Value = 250
$ws.Range("K140").Value = 12000
$ws.Range("L140").Value = 14000
$ws.Range("M140").Value = 12800
$ws.Range("N140").Value = "`$/malla 25 kilos"
$ws.Range("O140").Value = "Región Metropolitana"
$ws.Range("P140").Value = 512
$ws.Range("Q140").Value = 25
$ws.Range("R140").Value = "Hortaliza"

# New row 141 data
$ws.Range("A141").Value = 6
$ws.Range("B141").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C141").Value = "Metropolitana"
$ws.Range("D141").Value = 44505
$ws.Range("E141").Value = 13
$ws.Range("F141").Value = 100112022
$ws.Range("G141").Value = "Arveja Verde"
$ws.Range("H141").Value = "Sin especificar"
$ws.Range("I141").Value = "Primera"
$ws.Range("J141").Value = 250
$ws.Range("K141").Value = 10000
$ws.Range("L141").Value = 12000
$ws.Range("M141").Value = 11040
$ws.Range("N141").Value = "`$/saco 25 kilos"
$ws.Range("O141").Value = "Región Metropolitana"
$ws.Range("P141").Value = 442
$ws.Range("Q141").Value = 25
$ws.Range("R141").Value = "Hortaliza"
